$wb = $excel.ActiveWorkbook

$wsLessons   = $wb.Worksheets.Item("Lessons")
$wsQuestions = $wb.Worksheets.Item("Questions")
$wsAnswers   = $wb.Worksheets.Item("Answers")
$wsResults   = $wb.Worksheets.Item("Results")
$wsCes       = $wb.Worksheets.Item("Ces")

$dateTimeExample = "'2025-06-14 22:27:18.182458"
$hashExample     = "op8+4WOsllnwVsSOdtC7Jg=="

# --- Lessons: rows 8 & 9 ---
$wsLessons.Range("A8").Value = "date_time"
$wsLessons.Range("B8").Value = "datetime"
$wsLessons.Range("C8").Value = $false
$wsLessons.Range("D8").Value = $dateTimeExample

$wsLessons.Range("A9").Value = "date_time_hash"
$wsLessons.Range("B9").Value = "varchar(32)"
$wsLessons.Range("C9").Value = $false
$wsLessons.Range("D9").Value = $hashExample
$wsLessons.Range("E9").Value = "An MD5 hash of the date time value."

$wsLessons.Range("E8").Value = "The date time the lesson was created."

# --- Questions: rows 10 & 11 ---
$wsQuestions.Range("A10").Value = "date_time"
$wsQuestions.Range("B10").Value = "datetime"
$wsQuestions.Range("C10").Value = $false
$wsQuestions.Range("D10").Value = $dateTimeExample

$wsQuestions.Range("A11").Value = "date_time_hash"
$wsQuestions.Range("B11").Value = "varchar(32)"
$wsQuestions.Range("C11").Value = $false
$wsQuestions.Range("D11").Value = $hashExample
$wsQuestions.Range("E11").Value = "An MD5 hash of the date time value."

$wsQuestions.Range("E10").Value = "The date time the question was created."

# --- Answers: rows 7 & 8 ---
$wsAnswers.Range("A7").Value = "date_time"
$wsAnswers.Range("B7").Value = "datetime"
$wsAnswers.Range("C7").Value = $false
$wsAnswers.Range("D7").Value = $dateTimeExample

$wsAnswers.Range("A8").Value = "date_time_hash"
$wsAnswers.Range("B8").Value = "varchar(32)"
$wsAnswers.Range("C8").Value = $false
$wsAnswers.Range("D8").Value = $hashExample
$wsAnswers.Range("E8").Value = "An MD5 hash of the date time value."

$wsAnswers.Range("E7").Value = "The date time the answer was created."

# --- Results: rows 6 & 7 ---
$wsResults.Range("A6").Value = "date_time"
$wsResults.Range("B6").Value = "datetime"
$wsResults.Range("C6").Value = $false
$wsResults.Range("D6").Value = $dateTimeExample

$wsResults.Range("A7").Value = "date_time_hash"
$wsResults.Range("B7").Value = "varchar(32)"
$wsResults.Range("C7").Value = $false
$wsResults.Range("D7").Value = $hashExample
$wsResults.Range("E7").Value = "An MD5 hash of the date time value."

$wsResults.Range("E6").Value = "The date time the result was created."

# --- Final navigation / selection state to match the saved view ---
[void]$wsLessons.Activate()
$wsLessons.Rows(9).Select() | Out-Null

[void]$wsAnswers.Activate()
$wsAnswers.Rows(8).Select() | Out-Null

[void]$wsResults.Activate()
$wsResults.Range("E18").Select() | Out-Null

[void]$wsCes.Activate()
$wsCes.Range("A7").Select() | Out-Null

[void]$wsQuestions.Activate()
$wsQuestions.Range("D17").Select() | Out-Null
